$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $oldStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $oldStyle
}

# Row 8 and Row 9 swap: Cardano <-> Dogecoin (B, C, D columns), plus E updates
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D8" "0.0607"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.244"
$ws.Range("E9").Value = "  -0.44%  "

# Remaining price (D) and volume (E) updates
$ws.Range("D2").Value = "26.317.55"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.602.88"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "212.68"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.14%  "
Set-TextValue "D10" "18.97"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.827.48"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "1.601.42"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  -2.38%  "
Set-TextValue "D16" "63.62"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "26.305.07"
$ws.Range("E17").Value = "  +0.44%  "
Set-TextValue "D18" "227.74"
$ws.Range("E18").Value = "  +7.03%  "
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue "D20" "7.60"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("E21").Value = "  +0.09%  "
Set-TextValue "D22" "4.32"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  +0.75%  "
Set-TextValue "D24" "8.96"
$ws.Range("E24").Value = "  -0.15%  "
Set-TextValue "D25" "145.47"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  +1.27%  "
Set-TextValue "D29" "15.46"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "1.448.71"
$ws.Range("E33").Value = "  +7.86%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("E38").Value = "  -1.00%  "
Set-TextValue "D39" "0.825"
$ws.Range("E39").Value = "  +1.00%  "
Set-TextValue "D40" "5.83"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +1.99%  "
Set-TextValue "D43" "0.927"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("D44").Value = "1.740.94"
$ws.Range("E44").Value = "  +1.10%  "
Set-TextValue "D45" "0.758"
$ws.Range("E45").Value = "  -1.27%  "
Set-TextValue "D46" "60.96"
$ws.Range("E46").Value = "  -0.10%  "
Set-TextValue "D47" "87.67"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +0.24%  "
Set-TextValue "D49" "1.49"
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -3.45%  "
